# Update odds data in Sheet1 to reflect the latest FlashScore scrape
# (commit: "Atualizando o arquivo XLSX")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.48
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 7
$ws.Range("K2").Value = 2.2
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("W2").Value = 5.5
$ws.Range("Z2").Value = 9.5
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 81
$ws.Range("AM2").Value = 51
$ws.Range("AQ2").Value = 23
$ws.Range("AT2").Value = 2.63
$ws.Range("BA2").Value = 201

# --- Row 4 updates ---
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.47
$ws.Range("K4").Value = 2.07
$ws.Range("L4").Value = 3.05
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 2.9
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.75
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 15
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 5.9
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 450
$ws.Range("AI4").Value = 12.5
$ws.Range("AK4").Value = 28
$ws.Range("AM4").Value = 30
$ws.Range("AN4").Value = 4.75
$ws.Range("AP4").Value = 19.5
$ws.Range("AQ4").Value = 65
$ws.Range("AT4").Value = 2.6
$ws.Range("AY4").Value = 19.5
